# Regenerate save_data to use K (strikeouts) instead of Strike# and
# recompute/write the s_vals (G column) for each row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 1
    3  = 0
    4  = 2
    5  = 3
    6  = 0
    7  = 2
    8  = 1
    9  = 1
    10 = 5
    11 = 2
    12 = 4
    13 = 4
    14 = 7
    15 = 6
    16 = 9
    17 = 6
    18 = 9
    19 = 7
    20 = 8
    21 = 7
    22 = 8
    23 = 4
    24 = 4
    25 = 7
    26 = 5
    27 = 3
    28 = 6
    29 = 2
    30 = 2
    31 = 3
    32 = 4
    33 = 5
    34 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
